# CIERRE 20 DIC 23
# Adds the 19/12/2023 and 20/12/2023 movements to the ledger and extends
# the running-balance table (with blank placeholder rows) through row 83.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- 19 DIC 23: retiro de 6 botellones -----------------------------------
$ws.Range("B69").Value = 45279
$ws.Range("C69").Value = "6 botellones"
$ws.Range("D69").Value = -318
# E69 already holds the running-balance formula (=E68+D69); it recalculates
# automatically once D69 has a value.

# --- 20 DIC 23: pago ------------------------------------------------------
$ws.Range("B70").Value = 45280
$ws.Range("C70").Value = "a comprobar"
$ws.Range("D70").Value = 500
$ws.Range("E70").Formula = "=E69+D70"

# --- Extend the table with blank placeholder rows through row 83 ---------
# Row 71 through 83 get the same look & feel (number formats / borders) as
# row 70, copied across in one shot.
$ws.Range("B70:E70").Copy()
$ws.Range("B71:E83").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Running-balance formula for every new blank row except the very last one
# (row 83 stays without a formula, same as the old trailing row 70 used to).
$ws.Range("E71").Formula = "=E70+D71"
$ws.Range("E72").Formula = "=E71+D72"
$ws.Range("E73").Formula = "=E72+D73"
$ws.Range("E74").Formula = "=E73+D74"
$ws.Range("E75").Formula = "=E74+D75"
$ws.Range("E76").Formula = "=E75+D76"
$ws.Range("E77").Formula = "=E76+D77"
$ws.Range("E78").Formula = "=E77+D78"
$ws.Range("E79").Formula = "=E78+D79"
$ws.Range("E80").Formula = "=E79+D80"
$ws.Range("E81").Formula = "=E80+D81"
$ws.Range("E82").Formula = "=E81+D82"
$ws.Range("E83").ClearContents()

# --- View state: keep the frozen header and move the active cell down ----
$ws.Range("D71").Select()
